# The deck's design theme (ppt/theme/theme1.xml, bound to the slide master)
# is switched from the "Integral" / "Red Violet" palette over to the
# default PowerPoint "Office Theme" palette ("Office" colour scheme).
#
# PowerPoint's DrawingML colour scheme has 12 slots, in this fixed order:
#   dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
# `ThemeColorScheme` (reached off any Slide in the deck, since there is a
# single shared master/theme) exposes exactly those 12 slots and writes
# straight through to the theme part's <a:clrScheme>, so that's what we
# drive here.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target "Office" colour scheme values, dk1..folHlink (hex RRGGBB).
$officeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $hex = $officeColors[$i - 1]
    $rr = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $gg = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $bb = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # VBA/COM RGB values are packed 0x00BBGGRR.
    $tcs.Colors($i).RGB = ($bb * 65536) + ($gg * 256) + $rr
}

Write-Output "Updated theme colour scheme to Office defaults."
